# Apply updates to the "想去人数" (F) and "最低票价" (G) columns
# across the four sheets, matching the gh-pages data refresh commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 2357
$ws.Cells.Item(3, 6).Value = 532
$ws.Cells.Item(4, 6).Value = 209
$ws.Cells.Item(5, 6).Value = 355
$ws.Cells.Item(6, 6).Value = 355
$ws.Cells.Item(7, 6).Value = 578
$ws.Cells.Item(9, 6).Value = 795
$ws.Cells.Item(10, 6).Value = 541
$ws.Cells.Item(11, 6).Value = 816
$ws.Cells.Item(12, 6).Value = 389
$ws.Cells.Item(13, 6).Value = 99
$ws.Cells.Item(15, 6).Value = 22
$ws.Cells.Item(16, 6).Value = 1031
$ws.Cells.Item(17, 6).Value = 21181
$ws.Cells.Item(17, 7).Value = "已售罄"
$ws.Cells.Item(18, 6).Value = 822
$ws.Cells.Item(19, 6).Value = 77
$ws.Cells.Item(20, 6).Value = 262
$ws.Cells.Item(21, 6).Value = 303
$ws.Cells.Item(22, 6).Value = 179
$ws.Cells.Item(23, 6).Value = 166
$ws.Cells.Item(25, 6).Value = 16
$ws.Cells.Item(26, 6).Value = 230
$ws.Cells.Item(28, 6).Value = 356
$ws.Cells.Item(29, 6).Value = 155

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(6, 6).Value = 207
$ws.Cells.Item(7, 6).Value = 228
$ws.Cells.Item(8, 6).Value = 3438
$ws.Cells.Item(10, 6).Value = 101
$ws.Cells.Item(14, 6).Value = 124
$ws.Cells.Item(16, 6).Value = 3869

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 109
$ws.Cells.Item(4, 6).Value = 619
$ws.Cells.Item(5, 6).Value = 213

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 109
$ws.Cells.Item(5, 6).Value = 2357
$ws.Cells.Item(6, 6).Value = 619
$ws.Cells.Item(7, 6).Value = 532
$ws.Cells.Item(8, 6).Value = 209
$ws.Cells.Item(9, 6).Value = 355
$ws.Cells.Item(10, 6).Value = 355
$ws.Cells.Item(11, 6).Value = 578
$ws.Cells.Item(14, 6).Value = 3
$ws.Cells.Item(16, 6).Value = 207
$ws.Cells.Item(17, 6).Value = 213
$ws.Cells.Item(18, 6).Value = 795
$ws.Cells.Item(19, 6).Value = 541
$ws.Cells.Item(20, 6).Value = 816
$ws.Cells.Item(21, 6).Value = 389
$ws.Cells.Item(22, 6).Value = 99
$ws.Cells.Item(24, 6).Value = 22
$ws.Cells.Item(25, 6).Value = 1031
$ws.Cells.Item(26, 6).Value = 21182
$ws.Cells.Item(26, 7).Value = "已售罄"
$ws.Cells.Item(27, 6).Value = 228
$ws.Cells.Item(28, 6).Value = 3438
$ws.Cells.Item(30, 6).Value = 101
$ws.Cells.Item(32, 6).Value = 823
$ws.Cells.Item(33, 6).Value = 77
$ws.Cells.Item(34, 6).Value = 262
$ws.Cells.Item(37, 6).Value = 303
$ws.Cells.Item(38, 6).Value = 179
$ws.Cells.Item(39, 6).Value = 166
$ws.Cells.Item(41, 6).Value = 16
$ws.Cells.Item(42, 6).Value = 124
$ws.Cells.Item(44, 6).Value = 231
$ws.Cells.Item(46, 6).Value = 357
$ws.Cells.Item(47, 6).Value = 155
$ws.Cells.Item(48, 6).Value = 3870

